$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 382 (shifts old rows 382-388 down to 384-390)
$ws.Rows.Item(382).Insert()
$ws.Rows.Item(382).Insert()

# New row 382: Femacal de La Calera - Frutilla - Especial
$ws.Cells.Item(382, 1).Value = 3
$ws.Cells.Item(382, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(382, 3).Value = "Coquimbo"
$ws.Cells.Item(382, 4).Value = 44890
$ws.Cells.Item(382, 5).Value = 5
$ws.Cells.Item(382, 6).Value = "Fruta"
$ws.Cells.Item(382, 7).Value = 100101
$ws.Cells.Item(382, 8).Value = "Berries"
$ws.Cells.Item(382, 9).Value = 100112025
$ws.Cells.Item(382, 10).Value = "Frutilla"
$ws.Cells.Item(382, 11).Value = "Sin especificar"
$ws.Cells.Item(382, 12).Value = "Especial"
$ws.Cells.Item(382, 13).Value = 54
$ws.Cells.Item(382, 14).Value = 8000
$ws.Cells.Item(382, 15).Value = 8000
$ws.Cells.Item(382, 16).Value = 8000
$ws.Cells.Item(382, 17).Value = "`$/bandeja 7 kilos"
$ws.Cells.Item(382, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(382, 19).Value = 1143
$ws.Cells.Item(382, 20).Value = 7

# New row 383: Femacal de La Calera - Frutilla - Primera
$ws.Cells.Item(383, 1).Value = 3
$ws.Cells.Item(383, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(383, 3).Value = "Coquimbo"
$ws.Cells.Item(383, 4).Value = 44890
$ws.Cells.Item(383, 5).Value = 5
$ws.Cells.Item(383, 6).Value = "Fruta"
$ws.Cells.Item(383, 7).Value = 100101
$ws.Cells.Item(383, 8).Value = "Berries"
$ws.Cells.Item(383, 9).Value = 100112025
$ws.Cells.Item(383, 10).Value = "Frutilla"
$ws.Cells.Item(383, 11).Value = "Sin especificar"
$ws.Cells.Item(383, 12).Value = "Primera"
$ws.Cells.Item(383, 13).Value = 50
$ws.Cells.Item(383, 14).Value = 6000
$ws.Cells.Item(383, 15).Value = 6000
$ws.Cells.Item(383, 16).Value = 6000
$ws.Cells.Item(383, 17).Value = "`$/bandeja 7 kilos"
$ws.Cells.Item(383, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(383, 19).Value = 857
$ws.Cells.Item(383, 20).Value = 7
